# refactor: move code executions into main script
# Adds the results of a new AYTO episode: updates the matching_night_table
# with night 6 light counts, records a new perfect match (Aurelia+Josua)
# and several new no-matches, then leaves the "no_matches" sheet active.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: matching_night_table - new "night 6" data points
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("matching_night_table")

$ws1.Range("B3").Value = "4, 6, 7"
$ws1.Range("L4").Value = "4, 6, 7"
$ws1.Range("C5").Value = "3, 4, 5"
$ws1.Range("H6").Value = "5, 7"
$ws1.Range("B7").Value = 9
$ws1.Range("C7").Value = 0
$ws1.Range("H7").Value = "1, 3, 6"
$ws1.Range("K7").Value = "2, 4, 5, 7"
$ws1.Range("G11").Value = "1, 2, 4, 5, 6"

$ws1.Columns.Item(7).ColumnWidth = 9.1
$ws1.Columns.Item(10).ColumnWidth = 15.25

# ---------------------------------------------------------------------
# Sheet 3: perfect_matches - new perfect match found on night 9
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("perfect_matches")

$ws3.Range("A3").Value = 9
$ws3.Range("B3").Value = "Aurelia+Josua"
$ws3.Columns.Item(2).ColumnWidth = 14.92

# ---------------------------------------------------------------------
# Sheet 4: no_matches - updated/expanded no-match list
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("no_matches")

$ws4.Range("A2").Value = 1
$ws4.Range("B2").Value = "Finnja+Danilo"
$ws4.Range("A3").Value = 2
$ws4.Range("B3").Value = "Walentina+Tommy"
$ws4.Range("A4").Value = 8
$ws4.Range("B4").Value = "Melina+Tommy"
$ws4.Range("A5").Value = 5
$ws4.Range("B5").Value = "Finnja+Salvo"
$ws4.Range("A6").Value = 6
$ws4.Range("B6").Value = "Finnja+Eugen"
$ws4.Range("A7").Value = 7
$ws4.Range("B7").Value = "Steffi+Eugen"
$ws4.Range("A8").Value = 8
$ws4.Range("B8").Value = "Sarah+Josua"

$ws4.Columns.Item(2).ColumnWidth = 16.92

# ---------------------------------------------------------------------
# Selections: restore each sheet's cursor position, then leave
# "no_matches" as the active tab/sheet (done last so it "wins").
# ---------------------------------------------------------------------
[void]$ws1.Range("F15").Select()
[void]$ws3.Range("A3").Select()

[void]$ws4.Activate()
[void]$ws4.Range("E11").Select()
